$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 8 (this shifts all subsequent rows down by one,
# turning the old rows 8-12, 17-20, 25-29 into 9-13, 18-21, 26-30)
$ws.Rows.Item(8).Insert()

# Fill in the new "Dtype:" column (C) and a couple of additional "Info/comments" (F)
# and "Primary key" (D) cells describing the DB schema.

# Patient_info table
$ws.Range("C2").Value = "varchar?"
$ws.Range("D2").Value = "X"
$ws.Range("F2").Value = "e.g. 432_13"

$ws.Range("C3").Value = "Int"
$ws.Range("F3").Value = "e.g. 1995"

$ws.Range("C4").Value = "Varchar"

$ws.Range("C5").Value = "Int"

$ws.Range("F6").Value = "F/M"

$ws.Range("F7").Value = "Type of panel (e.g. APN, HSP, Filtex, Exome, ALS)"

# Seq_variants table
$ws.Range("C9").Value = "Varchar"
$ws.Range("C10").Value = "Int"
$ws.Range("C11").Value = "Int"

# Interpretations table
$ws.Range("C18").Value = "Int"
$ws.Range("C19").Value = "Varchar"
$ws.Range("C20").Value = "Date"
$ws.Range("C21").Value = "?"

# Runs table
$ws.Range("C26").Value = "Varchar"
$ws.Range("C27").Value = "Date"
$ws.Range("C28").Value = "Float"
$ws.Range("C29").Value = "Float"
$ws.Range("C30").Value = "Float"

# Update the active cell selection to match the final state
$ws.Range("F23").Select()
